# Update the raw metric inputs on the "Metrics" sheet. All the dependent
# formulas on other sheets (e.g. "today"!B11:B22, E11:E22, F11:F22) recompute
# automatically from these via the normal dependency graph.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metrics")

$ws1.Range("B2").Value  = 247708.03999999998
$ws1.Range("B3").Value  = 185007.76
$ws1.Range("B4").Value  = 63707.150000000009
$ws1.Range("B5").Value  = 10049
$ws1.Range("B6").Value  = 5883578.7699999986
$ws1.Range("B7").Value  = 4955725.3899999997
$ws1.Range("B8").Value  = 1727798.97
$ws1.Range("B9").Value  = 230326
$ws1.Range("B10").Value = 34348959.759999998
$ws1.Range("B11").Value = 32231000.550000001
$ws1.Range("B12").Value = 12009521.01
$ws1.Range("B13").Value = 1327956

# Move the selection/active-cell on "Metrics" to D17, then switch the active
# tab to "today" (its own selection, E9, is left untouched).
$ws1.Activate()
[void]$ws1.Range("D17").Select()

$ws2 = $wb.Worksheets.Item("today")
$ws2.Activate()
